$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalize inconsistent "District" values (column G) to the standard
# "Tumakuru (Tumkur)" label for the rows that still hold stray/old data.
$rows = @(5, 14, 20, 23, 38, 40, 42, 57, 59)
foreach ($r in $rows) {
    $ws.Range("G$r").Value = "Tumakuru (Tumkur)"
}
